# Actualiza base de datos EC: intercambia los valores de "Periodo Mora" y
# "Valor Mora" entre las filas 17 y 18 (mismo trabajador, dos periodos en mora).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 17 pasa a tener el periodo/valor que tenia la fila 18
$ws.Range("E17").Value = "2203"
$ws.Range("F17").Value = 35000

# Fila 18 pasa a tener el periodo/valor que tenia la fila 17
$ws.Range("E18").Value = "2204"
$ws.Range("F18").Value = 8400
